# ---------------------------------------------------------------------------
# Update EIA data to 2017
#
# Refreshes Table 6.2 "Coal Consumption by Sector" with the figures published
# in the June 2018 Monthly Energy Review (vs. the May 2017 release previously
# embedded in this workbook):
#   - masthead / release-date text on the "Monthly Data" sheet (shared by both
#     sheets via the shared-strings table)
#   - revised figures for the monthly rows that were still provisional in the
#     old release (Jan 2016 - Feb 2017)
#   - newly published monthly rows through March 2018
#   - the revised calendar-year-2016 annual total and the newly published
#     calendar-year-2017 annual row
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsMonthly = $wb.Worksheets.Item("Monthly Data")
$wsAnnual  = $wb.Worksheets.Item("Annual Data")

# --- Masthead text (same shared strings are rendered on both worksheets) ---
$wsMonthly.Range("A2").Value = "June 2018 Monthly Energy Review"
$wsMonthly.Range("A6").Value = "Release Date: June 26, 2018"
$wsMonthly.Range("A7").Value = "Next Update: July 26, 2018"
$wsAnnual.Range("A2").Value = "June 2018 Monthly Energy Review"
$wsAnnual.Range("A6").Value = "Release Date: June 26, 2018"
$wsAnnual.Range("A7").Value = "Next Update: July 26, 2018"

# --- Monthly Data: revised figures for existing rows 529-542 (Jan 2016-Feb 2017) ---
$monthlyRevisions = @(
    @(529, "C", 75.262, "D", 74.912, "E", 150.174, "G", 1397.261, "H", 1651.675, "I", 3048.936, "J", 4377.419, "L", 62134.631, "M", 66662.224),
    @(530, "C", 75.387, "D", 75.036, "E", 150.423, "G", 1282.036, "H", 1755.363, "I", 3037.399, "J", 4398.844, "L", 50661.45, "M", 55210.717),
    @(531, "C", 74.003, "D", 73.658, "E", 147.661, "G", 1274.833, "H", 1770.31, "I", 3045.143, "J", 4478.8, "L", 39948.145, "M", 44574.606),
    @(532, "C", 45.553, "D", 28.657, "E", 74.21, "G", 1075.583, "H", 1750.917, "I", 2826.5, "J", 4150.531, "L", 39158.963, "M", 43383.704),
    @(533, "C", 36.542, "D", 22.989, "E", 59.531, "G", 1178.029, "H", 1656.568, "I", 2834.597, "J", 4201.467, "L", 45081.935, "M", 49342.933),
    @(534, "C", 46.166, "D", 29.043, "E", 75.209, "G", 1242.702, "H", 1578.086, "I", 2820.788, "J", 4225.606, "L", 63250.414, "M", 67551.229),
    @(535, "C", 46.143, "D", 17.383, "E", 63.526, "G", 1321.161, "H", 1514.584, "I", 2835.745, "J", 4268.285, "L", 74236.728, "M", 78568.539),
    @(536, "C", 49.413, "D", 18.615, "E", 68.028, "G", 1292.162, "H", 1529.738, "I", 2821.9, "J", 4216.578, "L", 73889.93, "M", 78174.536),
    @(537, "C", 49.607, "D", 18.688, "E", 68.295, "G", 1157.257, "H", 1668.331, "I", 2825.588, "J", 4161.387, "L", 62385.216, "M", 66614.898),
    @(538, "C", 49.856, "D", 37.991, "E", 87.847, "G", 1126.246, "H", 1782.495, "I", 2908.741, "J", 4243.411, "L", 54621.445, "M", 58952.703),
    @(539, "C", 59.538, "D", 45.368, "E", 104.906, "G", 1092.706, "H", 1830.459, "I", 2923.165, "J", 4249.133, "L", 48179.203, "M", 52533.242),
    @(540, "C", 75.425, "D", 57.474, "E", 132.899, "G", 1280.244, "H", 1640.041, "I", 2920.285, "J", 4362.034, "L", 65006.425, "M", 69501.358),
    @(541, "C", 66.184, "D", 72.161, "E", 138.345, "F", 1430.645, "G", 1290.391, "H", 1553.562, "I", 2843.953, "J", 4274.598, "L", 63547.714, "M", 67960.657),
    @(542, "C", 53.674, "D", 58.521, "E", 112.195, "F", 1367.727, "G", 1087.427, "H", 1766.993, "I", 2854.42, "J", 4222.147, "L", 47964.848, "M", 52299.19)
)

foreach ($entry in $monthlyRevisions) {
    $r = $entry[0]
    for ($i = 1; $i -lt $entry.Count; $i += 2) {
        $colLetter = $entry[$i]
        $value = $entry[$i + 1]
        $wsMonthly.Range("$colLetter$r").Value = $value
    }
}

# --- Monthly Data: newly published rows 543-555 (Mar 2017-Mar 2018) ---
# Each new row inherits formatting (incl. the "yyyy mmmm" date style on column A)
# from the row immediately above it, then the values are overwritten.
$monthlyNewRows = @(
    @(543, 42795, "Not Available", 58.423, 63.699, 122.122, 1437.669, 1172.172, 1664.293, 2836.465, 4274.134, 0, 48825.958, 53222.214),
    @(544, 42826, "Not Available", 40.473, 24.9, 65.373, 1440.81, 1067.605, 1629.855, 2697.46, 4138.27, 0, 44323.847, 48527.49),
    @(545, 42856, "Not Available", 39.962, 24.586, 64.548, 1482.486, 1098.283, 1604.724, 2703.007, 4185.493, 0, 50926.005, 55176.046),
    @(546, 42887, "Not Available", 45.507, 27.997, 73.504, 1401.664, 1094.108, 1617.146, 2711.254, 4112.918, 0, 58951.924, 63138.346),
    @(547, 42917, "Not Available", 53.309, 16.546, 69.855, 1494.46, 1047.123, 1838.301, 2885.424, 4379.884, 0, 69900.111, 74349.85),
    @(548, 42948, "Not Available", 48.549, 15.069, 63.618, 1528.056, 1064.994, 1807.254, 2872.248, 4400.304, 0, 65933.994, 70397.916),
    @(549, 42979, "Not Available", 47.069, 14.609, 61.678, 1468.767, 1030.015, 1809.249, 2839.264, 4308.031, 0, 54779.784, 59149.493),
    @(550, 43009, "Not Available", 42.669, 37.703, 80.372, 1469.57, 1149.103, 1641.732, 2790.835, 4260.405, 0, 50214.467, 54555.244),
    @(551, 43040, "Not Available", 49.578, 43.808, 93.386, 1456.863, 1142.217, 1650.271, 2792.488, 4249.351, 0, 50992.13, 55334.867),
    @(552, 43070, "Not Available", 61.661, 54.486, 116.147, 1558.946, 1180.736, 1605.369, 2786.105, 4345.051, 0, 58388.345, 62849.543),
    @(553, 43101, "Not Available", 68.861, 35.059, 103.92, 1689.106, 1269.77, 1667.075, 2936.845, 4625.951, 0, 64650.176, 69380.047),
    @(554, 43132, "Not Available", 53.112, 50.81, 103.922, 1388.187, 1131.822, 1729.148, 2860.97, 4249.157, 0, 45823.067, 50176.146),
    @(555, 43160, "Not Available", 50.608, 7.14, 57.748, 1113.576, 1169.107, 1610.478, 2779.585, 3893.161, 0, 44495.503, 48446.411)
)

foreach ($entry in $monthlyNewRows) {
    $r = $entry[0]
    $prev = $r - 1
    $wsMonthly.Range("A${prev}:M${prev}").Copy()
    $wsMonthly.Range("A${r}:M${r}").PasteSpecial(-4122)
    $wsMonthly.Range("A$r").Value = $entry[1]
    $wsMonthly.Range("B$r").Value = $entry[2]
    $wsMonthly.Range("C$r").Value = $entry[3]
    $wsMonthly.Range("D$r").Value = $entry[4]
    $wsMonthly.Range("E$r").Value = $entry[5]
    $wsMonthly.Range("F$r").Value = $entry[6]
    $wsMonthly.Range("G$r").Value = $entry[7]
    $wsMonthly.Range("H$r").Value = $entry[8]
    $wsMonthly.Range("I$r").Value = $entry[9]
    $wsMonthly.Range("J$r").Value = $entry[10]
    $wsMonthly.Range("K$r").Value = $entry[11]
    $wsMonthly.Range("L$r").Value = $entry[12]
    $wsMonthly.Range("M$r").Value = $entry[13]
}

# --- Annual Data: revised calendar-year-2016 total (row 80) ---
$annualRevision = @("C", 682.895, "D", 499.814, "E", 1182.709, "G", 14720.22, "H", 20128.567, "I", 34848.787, "J", 51333.495, "L", 678554.486, "M", 731070.69)
for ($i = 0; $i -lt $annualRevision.Count; $i += 2) {
    $colLetter = $annualRevision[$i]
    $value = $annualRevision[$i + 1]
    $wsAnnual.Range("${colLetter}80").Value = $value
}

# --- Annual Data: newly published calendar-year-2017 row (81) ---
$annualNewRow = @(2017, "Not Available", 607.058, 454.085, 1061.143, 17537.663, 13424.174, 20188.749, 33612.923, 51150.586, 0, 664749.129, 716960.858)
$wsAnnual.Range("A79:M79").Copy()
$wsAnnual.Range("A81:M81").PasteSpecial(-4122)
$wsAnnual.Range("A81").Value = $annualNewRow[0]
$wsAnnual.Range("B81").Value = $annualNewRow[1]
$wsAnnual.Range("C81").Value = $annualNewRow[2]
$wsAnnual.Range("D81").Value = $annualNewRow[3]
$wsAnnual.Range("E81").Value = $annualNewRow[4]
$wsAnnual.Range("F81").Value = $annualNewRow[5]
$wsAnnual.Range("G81").Value = $annualNewRow[6]
$wsAnnual.Range("H81").Value = $annualNewRow[7]
$wsAnnual.Range("I81").Value = $annualNewRow[8]
$wsAnnual.Range("J81").Value = $annualNewRow[9]
$wsAnnual.Range("K81").Value = $annualNewRow[10]
$wsAnnual.Range("L81").Value = $annualNewRow[11]
$wsAnnual.Range("M81").Value = $annualNewRow[12]
